$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.874.61"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "3.852.23"
$ws.Range("E3").Value = "  +1.66%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "697.14"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.65"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").Value = "3.851.23"
$ws.Range("E7").Value = "  +1.67%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("E10").Value = "  -0.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.22"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.03%  "
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000257"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.26"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.40%  "
$ws.Range("D15").Value = "4.498.44"
$ws.Range("E15").Value = "  +1.61%  "
$ws.Range("D16").Value = "3.891.45"
$ws.Range("E16").Value = "  +2.67%  "
$ws.Range("D17").Value = "70.915.44"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.39"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.73%  "
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "498.07"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.65"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.96%  "
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.85"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000147"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.63"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.22"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.13"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.76%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.15"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.09%  "
$ws.Range("B30").Value = "Dai"
$ws.Range("C30").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.27"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.54"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  +2.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.19"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").Value = "3.806.12"
$ws.Range("E36").Value = "  +1.78%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.37"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +6.40%  "
$ws.Range("E40").Value = "  +8.69%  "
$ws.Range("E41").Value = "  -2.36%  "
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "163.99"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.88%  "
$ws.Range("B46").Value = "FLOKI"
$ws.Range("C46").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000310"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -5.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "49.04"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.39"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.24%  "
$ws.Range("B49").Value = "TheGraph"
$ws.Range("C49").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.300"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("B50").Value = "Arweave"
$ws.Range("C50").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "43.42"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.99%  "
$ws.Range("E51").Value = "  +1.35%  "
